$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 3
    4  = 5
    5  = 5
    6  = 4
    7  = 6
    8  = 3
    9  = 5
    10 = 3
    11 = 4
    12 = 4
    13 = 5
    14 = 7
    15 = 3
    16 = 2
    17 = 8
    18 = 4
    19 = 10
    20 = 2
    21 = 5
    22 = 4
    23 = 6
    24 = 10
    25 = 6
    26 = 2
    27 = 6
    28 = 4
    29 = 3
    30 = 3
    31 = 6
    32 = 5
    33 = 4
    34 = 4
    35 = 3
    36 = 5
    37 = 5
    38 = 2
    39 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
